$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value2 = 1445.6428
$ws.Range("I38").Value2 = 100.63636
$ws.Range("J38").Value2 = 6377.3335
$ws.Range("K38").Value2 = 301.90908
$ws.Range("L38").Value2 = 19132.0005
$ws.Range("M38").Value2 = 70.09091999999998
$ws.Range("N38").Value2 = -19876.0005
$ws.Range("H40").Value2 = 7666.3335
$ws.Range("J40").Value2 = 9999
$ws.Range("L40").Value2 = 9999
$ws.Range("N40").Value2 = -10349
$ws.Range("H41").Value2 = 698.8823
$ws.Range("I41").Value2 = 282.25
$ws.Range("K41").Value2 = 282.25
$ws.Range("M41").Value2 = 157.75
$ws.Range("H42").Value2 = 134.66667
$ws.Range("I42").Value2 = 121.6
$ws.Range("J42").Value2 = 200
$ws.Range("K42").Value2 = 364.8
$ws.Range("L42").Value2 = 600
$ws.Range("M42").Value2 = -134.8
$ws.Range("N42").Value2 = -1060
$ws.Range("H49").Value2 = 0
$ws.Range("I49").Value2 = 0
$ws.Range("K49").Value2 = 0
$ws.Range("M49").ClearContents()
$ws.Range("H70").Value2 = 4967.4644
$ws.Range("I70").Value2 = 5372.5
$ws.Range("J70").Value2 = 4936.3076
$ws.Range("K70").Value2 = 16117.5
$ws.Range("L70").Value2 = 14808.9228
$ws.Range("M70").Value2 = -15847.5
$ws.Range("N70").Value2 = -15348.9228
$ws.Range("H73").Value2 = 4967.4644
$ws.Range("I73").Value2 = 5372.5
$ws.Range("J73").Value2 = 4936.3076
$ws.Range("K73").Value2 = 16117.5
$ws.Range("L73").Value2 = 14808.9228
$ws.Range("M73").Value2 = -15181.5
$ws.Range("N73").Value2 = -16680.9228
$ws.Range("H99").Value2 = 100007000
$ws.Range("J99").Value2 = 250015870
$ws.Range("L99").Value2 = 750047610
$ws.Range("N99").Value2 = -750050606
$ws.Range("H101").Value2 = 1498
$ws.Range("I101").Value2 = 1498
$ws.Range("J101").Value2 = 0
$ws.Range("K101").Value2 = 4494
$ws.Range("L101").Value2 = 0
$ws.Range("M101").Value2 = -2872
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value2 = 350458.06
$ws.Range("I132").Value2 = 469761.84
$ws.Range("J132").Value2 = 14238.363
$ws.Range("K132").Value2 = 1409285.52
$ws.Range("L132").Value2 = 42715.089
$ws.Range("M132").Value2 = -1406755.52
$ws.Range("N132").Value2 = -47775.089
$ws.Range("H135").Value2 = 4109.1797
$ws.Range("I135").Value2 = 1823.5
$ws.Range("K135").Value2 = 16411.5
$ws.Range("M135").Value2 = -13876.5
$ws.Range("H138").Value2 = 5118.5415
$ws.Range("I138").Value2 = 4312.615
$ws.Range("J138").Value2 = 5296.1187
$ws.Range("K138").Value2 = 12937.845
$ws.Range("L138").Value2 = 15888.3561
$ws.Range("M138").Value2 = -7797.844999999999
$ws.Range("N138").Value2 = -26168.3561

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 16663.96
$ws.Range("I32").Value2 = 8277.151
$ws.Range("J32").Value2 = 32944.234
$ws.Range("K32").Value2 = 8277.151
$ws.Range("L32").Value2 = 32944.234
$ws.Range("M32").Value2 = -7990.151
$ws.Range("N32").Value2 = -33518.234
$ws.Range("H61").Value2 = 6290.3335
$ws.Range("I61").Value2 = 6296.143
$ws.Range("J61").Value2 = 6249.6665
$ws.Range("K61").Value2 = 6296.143
$ws.Range("L61").Value2 = 6249.6665
$ws.Range("M61").Value2 = -6084.143
$ws.Range("N61").Value2 = -6673.6665
$ws.Range("H74").Value2 = 4871
$ws.Range("I74").Value2 = 5709.1665
$ws.Range("J74").Value2 = 3434.1428
$ws.Range("K74").Value2 = 5709.1665
$ws.Range("L74").Value2 = 3434.1428
$ws.Range("M74").Value2 = -4835.1665
$ws.Range("N74").Value2 = -5182.1428
$ws.Range("H77").Value2 = 4871
$ws.Range("I77").Value2 = 5709.1665
$ws.Range("J77").Value2 = 3434.1428
$ws.Range("K77").Value2 = 28545.8325
$ws.Range("L77").Value2 = 17170.714
$ws.Range("M77").Value2 = -24177.8325
$ws.Range("N77").Value2 = -25906.714
$ws.Range("H97").Value2 = 3135.8
$ws.Range("I97").Value2 = 1441.4166
$ws.Range("K97").Value2 = 1441.4166
$ws.Range("M97").Value2 = -945.4166
$ws.Range("H132").Value2 = 598090.75
$ws.Range("I132").Value2 = 829300.25
$ws.Range("K132").Value2 = 2487900.75
$ws.Range("M132").Value2 = -2485370.75
$ws.Range("H136").Value2 = 6290.3335
$ws.Range("I136").Value2 = 6296.143
$ws.Range("J136").Value2 = 6249.6665
$ws.Range("K136").Value2 = 18888.429
$ws.Range("L136").Value2 = 18748.9995
$ws.Range("M136").Value2 = -16338.429
$ws.Range("N136").Value2 = -23848.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 10973.394
$ws.Range("I99").Value2 = 11572.647
$ws.Range("J99").Value2 = 10218.777
$ws.Range("K99").Value2 = 11572.647
$ws.Range("L99").Value2 = 10218.777
$ws.Range("M99").Value2 = -10074.647
$ws.Range("N99").Value2 = -13214.777
$ws.Range("H137").Value2 = 66999
$ws.Range("J137").Value2 = 66999
$ws.Range("L137").Value2 = 66999
$ws.Range("N137").Value2 = -77199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 10037.23
$ws.Range("I31").Value2 = 1899
$ws.Range("J31").Value2 = 15123.625
$ws.Range("K31").Value2 = 1899
$ws.Range("L31").Value2 = 15123.625
$ws.Range("M31").Value2 = -1604
$ws.Range("N31").Value2 = -15713.625
$ws.Range("H34").Value2 = 10037.23
$ws.Range("I34").Value2 = 1899
$ws.Range("J34").Value2 = 15123.625
$ws.Range("K34").Value2 = 1899
$ws.Range("L34").Value2 = 15123.625
$ws.Range("M34").Value2 = -1697
$ws.Range("N34").Value2 = -15527.625
$ws.Range("H62").Value2 = 9688.182000000001
$ws.Range("I62").Value2 = 9572
$ws.Range("J62").Value2 = 9998
$ws.Range("K62").Value2 = 9572
$ws.Range("L62").Value2 = 9998
$ws.Range("M62").Value2 = -8948
$ws.Range("N62").Value2 = -11246
$ws.Range("H65").Value2 = 9688.182000000001
$ws.Range("I65").Value2 = 9572
$ws.Range("J65").Value2 = 9998
$ws.Range("K65").Value2 = 47860
$ws.Range("L65").Value2 = 49990
$ws.Range("M65").Value2 = -44740
$ws.Range("N65").Value2 = -56230
$ws.Range("H105").Value2 = 50001708
$ws.Range("I105").Value2 = 50001708
$ws.Range("K105").Value2 = 50001708
$ws.Range("M105").Value2 = -49999961
$ws.Range("H132").Value2 = 9015.352999999999
$ws.Range("I132").Value2 = 7814.5557
$ws.Range("K132").Value2 = 23443.6671
$ws.Range("M132").Value2 = -20913.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 1437.129
$ws.Range("I5").Value2 = 763.8333
$ws.Range("K5").Value2 = 2291.4999
$ws.Range("M5").Value2 = -2179.4999
$ws.Range("H17").Value2 = 8021.5713
$ws.Range("I17").Value2 = 9050
$ws.Range("J17").Value2 = 7610.2
$ws.Range("K17").Value2 = 27150
$ws.Range("L17").Value2 = 22830.6
$ws.Range("M17").Value2 = -26981
$ws.Range("N17").Value2 = -23168.6
$ws.Range("H23").Value2 = 279.55554
$ws.Range("J23").Value2 = 288.2857
$ws.Range("L23").Value2 = 864.8571000000001
$ws.Range("N23").Value2 = -1334.8571
$ws.Range("H37").Value2 = 123218
$ws.Range("J37").Value2 = 123218
$ws.Range("L37").Value2 = 369654
$ws.Range("N37").Value2 = -369878
$ws.Range("H113").Value2 = 1977
$ws.Range("I113").Value2 = 1509.5555
$ws.Range("J113").Value2 = 2187.35
$ws.Range("K113").Value2 = 4528.666499999999
$ws.Range("L113").Value2 = 6562.049999999999
$ws.Range("M113").Value2 = -2358.666499999999
$ws.Range("N113").Value2 = -10902.05
$ws.Range("H135").Value2 = 1437.129
$ws.Range("I135").Value2 = 763.8333
$ws.Range("K135").Value2 = 6874.4997
$ws.Range("M135").Value2 = -4339.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 50009280
$ws.Range("J80").Value2 = 14709.75
$ws.Range("L80").Value2 = 14709.75
$ws.Range("N80").Value2 = -16705.75
$ws.Range("H83").Value2 = 50009280
$ws.Range("J83").Value2 = 14709.75
$ws.Range("L83").Value2 = 73548.75
$ws.Range("N83").Value2 = -83532.75
$ws.Range("H102").Value2 = 4434.278
$ws.Range("I102").Value2 = 3360.8865
$ws.Range("J102").Value2 = 6121.0356
$ws.Range("K102").Value2 = 3360.8865
$ws.Range("L102").Value2 = 6121.0356
$ws.Range("M102").Value2 = -1738.8865
$ws.Range("N102").Value2 = -9365.035599999999
$ws.Range("H127").Value2 = 0
$ws.Range("J127").Value2 = 0
$ws.Range("L127").Value2 = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value2 = 5137.5557
$ws.Range("I93").Value2 = 3920.4375
$ws.Range("J93").Value2 = 14874.5
$ws.Range("K93").Value2 = 3920.4375
$ws.Range("L93").Value2 = 14874.5
$ws.Range("M93").Value2 = -2672.4375
$ws.Range("N93").Value2 = -17370.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 10031.956
$ws.Range("J62").Value2 = 9756.056
$ws.Range("L62").Value2 = 9756.056
$ws.Range("N62").Value2 = -11004.056
$ws.Range("H65").Value2 = 10031.956
$ws.Range("J65").Value2 = 9756.056
$ws.Range("L65").Value2 = 48780.28
$ws.Range("N65").Value2 = -55020.28
$ws.Range("H122").Value2 = 9421.1
$ws.Range("I122").Value2 = 3586.16
$ws.Range("J122").Value2 = 19146
$ws.Range("K122").Value2 = 10758.48
$ws.Range("L122").Value2 = 57438
$ws.Range("M122").Value2 = -8308.48
$ws.Range("N122").Value2 = -62338
